# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The value 45181 (2023-09-12) is bumped by one day to 45182 (2023-09-13)
# for every row in the worksheet's used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
